$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 100001
$ws.Range("I12").Value = 100001
$ws.Range("K12").Value = 100001
$ws.Range("M12").Value = -99831

$ws.Range("H15").Value = 217.07
$ws.Range("I15").Value = 217.07
$ws.Range("K15").Value = 651.21
$ws.Range("M15").Value = -482.21

$ws.Range("H88").Value = 905.4583
$ws.Range("I88").Value = 876.3333
$ws.Range("J88").Value = 915.1667
$ws.Range("K88").Value = 876.3333
$ws.Range("L88").Value = 915.1667
$ws.Range("M88").Value = -470.3333
$ws.Range("N88").Value = -1727.1667

$ws.Range("H91").Value = 905.4583
$ws.Range("I91").Value = 876.3333
$ws.Range("J91").Value = 915.1667
$ws.Range("K91").Value = 876.3333
$ws.Range("L91").Value = 915.1667
$ws.Range("M91").Value = 527.6667
$ws.Range("N91").Value = -3723.1667

$ws.Range("H107").Value = 443.35715
$ws.Range("I107").Value = 422.74075
$ws.Range("K107").Value = 422.74075
$ws.Range("M107").Value = 1497.25925

$ws.Range("H124").Value = 79800
$ws.Range("J124").Value = 79800
$ws.Range("L124").Value = 79800
$ws.Range("N124").Value = -89620

$ws.Range("H130").Value = 79800
$ws.Range("J130").Value = 79800
$ws.Range("L130").Value = 79800
$ws.Range("N130").Value = -89840

$ws.Range("H138").Value = 3417.719
$ws.Range("J138").Value = 4146.4756
$ws.Range("L138").Value = 12439.4268
$ws.Range("N138").Value = -22719.4268

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7563.183
$ws.Range("I32").Value = 6961.9136
$ws.Range("J32").Value = 25000
$ws.Range("K32").Value = 6961.9136
$ws.Range("L32").Value = 25000
$ws.Range("M32").Value = -6674.9136
$ws.Range("N32").Value = -25574

$ws.Range("H45").Value = 1552.8182
$ws.Range("I45").Value = 1540.2307
$ws.Range("K45").Value = 1540.2307
$ws.Range("M45").Value = -1163.2307

$ws.Range("H56").Value = 32400
$ws.Range("J56").Value = 32400
$ws.Range("L56").Value = 32400
$ws.Range("N56").Value = -33884

$ws.Range("H110").Value = 1592.3077
$ws.Range("I110").Value = 1410
$ws.Range("J110").Value = 2200
$ws.Range("K110").Value = 1410
$ws.Range("L110").Value = 2200
$ws.Range("M110").Value = 635
$ws.Range("N110").Value = -6290

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1742.5278
$ws.Range("I94").Value = 1660.92
$ws.Range("J94").Value = 1928
$ws.Range("K94").Value = 1660.92
$ws.Range("L94").Value = 1928
$ws.Range("M94").Value = -1209.92
$ws.Range("N94").Value = -2830

$ws.Range("H107").Value = 3250.9167
$ws.Range("I107").Value = 3001.375
$ws.Range("J107").Value = 3750
$ws.Range("K107").Value = 3001.375
$ws.Range("L107").Value = 3750
$ws.Range("M107").Value = -1081.375
$ws.Range("N107").Value = -7590

$ws.Range("H134").Value = 48619.137
$ws.Range("I134").Value = 3423.6875
$ws.Range("J134").Value = 169140.33
$ws.Range("K134").Value = 10271.0625
$ws.Range("L134").Value = 507420.99
$ws.Range("M134").Value = -7736.0625
$ws.Range("N134").Value = -512490.99

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2186.2185
$ws.Range("I31").Value = 1556.0892
$ws.Range("J31").Value = 3324.516
$ws.Range("K31").Value = 1556.0892
$ws.Range("L31").Value = 3324.516
$ws.Range("M31").Value = -1261.0892
$ws.Range("N31").Value = -3914.516

$ws.Range("H34").Value = 2186.2185
$ws.Range("I34").Value = 1556.0892
$ws.Range("J34").Value = 3324.516
$ws.Range("K34").Value = 1556.0892
$ws.Range("L34").Value = 3324.516
$ws.Range("M34").Value = -1354.0892
$ws.Range("N34").Value = -3728.516

$ws.Range("H132").Value = 1999.7234
$ws.Range("I132").Value = 1441.2059
$ws.Range("J132").Value = 3460.4614
$ws.Range("K132").Value = 4323.6177
$ws.Range("L132").Value = 10381.3842
$ws.Range("M132").Value = -1793.6177
$ws.Range("N132").Value = -15441.3842

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 36098.54
$ws.Range("I14").Value = 36098.54
$ws.Range("K14").Value = 108295.62
$ws.Range("M14").Value = -108122.62

$ws.Range("H44").Value = 18650
$ws.Range("I44").Value = 18650
$ws.Range("K44").Value = 55950
$ws.Range("M44").Value = -55552

$ws.Range("H121").Value = 1212.5454
$ws.Range("I121").Value = 515.7143
$ws.Range("J121").Value = 1537.7333
$ws.Range("K121").Value = 1547.1429
$ws.Range("L121").Value = 4613.199900000001
$ws.Range("M121").Value = -237.1428999999998
$ws.Range("N121").Value = -7233.199900000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 292.57144
$ws.Range("I2").Value = 112
$ws.Range("J2").Value = 533.3333
$ws.Range("K2").Value = 112
$ws.Range("L2").Value = 533.3333
$ws.Range("M2").Value = 1
$ws.Range("N2").Value = -759.3333

$ws.Range("H97").Value = 1516.6296
$ws.Range("I97").Value = 1244.0454
$ws.Range("J97").Value = 2716
$ws.Range("K97").Value = 1244.0454
$ws.Range("L97").Value = 2716
$ws.Range("M97").Value = -748.0454
$ws.Range("N97").Value = -3708

$ws.Range("H126").Value = 3048.5715
$ws.Range("I126").Value = 1999.7778
$ws.Range("J126").Value = 3835.1667
$ws.Range("K126").Value = 5999.3334
$ws.Range("L126").Value = 11505.5001
$ws.Range("M126").Value = -3529.3334
$ws.Range("N126").Value = -16445.5001

$ws.Range("H132").Value = 63843.945
$ws.Range("I132").Value = 254294.75
$ws.Range("J132").Value = 9429.429
$ws.Range("K132").Value = 762884.25
$ws.Range("L132").Value = 28288.287
$ws.Range("M132").Value = -760354.25
$ws.Range("N132").Value = -33348.287

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3633.9285
$ws.Range("I7").Value = 3250
$ws.Range("J7").Value = 4145.8335
$ws.Range("K7").Value = 3250
$ws.Range("L7").Value = 4145.8335
$ws.Range("M7").Value = -3138
$ws.Range("N7").Value = -4369.8335

$ws.Range("H55").Value = 160632.48
$ws.Range("I55").Value = 267302.8
$ws.Range("J55").Value = 627
$ws.Range("K55").Value = 267302.8
$ws.Range("L55").Value = 627
$ws.Range("M55").Value = -267129.8
$ws.Range("N55").Value = -973

$ws.Range("H126").Value = 3633.9285
$ws.Range("I126").Value = 3250
$ws.Range("J126").Value = 4145.8335
$ws.Range("K126").Value = 9750
$ws.Range("L126").Value = 12437.5005
$ws.Range("M126").Value = -7280
$ws.Range("N126").Value = -17377.5005

$ws.Range("H132").Value = 4024.3
$ws.Range("I132").Value = 3754.8572
$ws.Range("J132").Value = 4653
$ws.Range("K132").Value = 11264.5716
$ws.Range("L132").Value = 13959
$ws.Range("M132").Value = -8734.5716
$ws.Range("N132").Value = -19019

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5433.6
$ws.Range("I122").Value = 3017.4443
$ws.Range("J122").Value = 7410.4546
$ws.Range("K122").Value = 9052.332900000001
$ws.Range("L122").Value = 22231.3638
$ws.Range("M122").Value = -6602.332900000001
$ws.Range("N122").Value = -27131.3638

$ws.Range("H132").Value = 3812.5483
$ws.Range("I132").Value = 3186.4285
$ws.Range("K132").Value = 9559.2855
$ws.Range("M132").Value = -7029.2855
